$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 1.02
$ws.Cells.Item(2, 3).Value = 1.01429419160473
$ws.Cells.Item(2, 4).Value = 1.020829187621129
$ws.Cells.Item(2, 5).Value = 1.016107839280098
$ws.Cells.Item(2, 9).Value = 1.025876530047811
$ws.Cells.Item(2, 10).Value = 1.019526252390308
$ws.Cells.Item(2, 11).Value = 1.023668316833596
$ws.Cells.Item(2, 12).Value = 1.01896099625379
$ws.Cells.Item(2, 14).Value = 1.010660548212845

$ws.Cells.Item(3, 2).Value = 1.02
$ws.Cells.Item(3, 3).Value = 1.015102552581614
$ws.Cells.Item(3, 4).Value = 1.0214192022894
$ws.Cells.Item(3, 5).Value = 1.016788647862587
$ws.Cells.Item(3, 9).Value = 1.025978726153472
$ws.Cells.Item(3, 10).Value = 1.019969873490211
$ws.Cells.Item(3, 11).Value = 1.024065247508632
$ws.Cells.Item(3, 12).Value = 1.019447457670659
$ws.Cells.Item(3, 14).Value = 1.010806075101252

$ws.Cells.Item(4, 2).Value = 1.02
$ws.Cells.Item(4, 3).Value = 1.015625708554958
$ws.Cells.Item(4, 4).Value = 1.021800415157856
$ws.Cells.Item(4, 5).Value = 1.017229679440857
$ws.Cells.Item(4, 9).Value = 1.026042775220996
$ws.Cells.Item(4, 10).Value = 1.020256383972522
$ws.Cells.Item(4, 11).Value = 1.024320827608154
$ws.Cells.Item(4, 12).Value = 1.01976205789824
$ws.Cells.Item(4, 14).Value = 1.010900060312771

$ws.Cells.Item(5, 2).Value = 1.02
$ws.Cells.Item(5, 3).Value = 1.015845663729786
$ws.Cells.Item(5, 4).Value = 1.021960539451098
$ws.Cells.Item(5, 5).Value = 1.017415207618748
$ws.Cells.Item(5, 9).Value = 1.026069202749065
$ws.Cells.Item(5, 10).Value = 1.020376701814905
$ws.Cells.Item(5, 11).Value = 1.024427970050174
$ws.Cells.Item(5, 12).Value = 1.019894272990309
$ws.Cells.Item(5, 14).Value = 1.010939528035934

$ws.Cells.Item(6, 2).Value = 1.02
$ws.Cells.Item(6, 3).Value = 1.015882596318142
$ws.Cells.Item(6, 4).Value = 1.021987416895099
$ws.Cells.Item(6, 5).Value = 1.017446365509986
$ws.Cells.Item(6, 9).Value = 1.026073610772076
$ws.Cells.Item(6, 10).Value = 1.020396895977633
$ws.Cells.Item(6, 11).Value = 1.024445941896196
$ws.Cells.Item(6, 12).Value = 1.019916469930888
$ws.Cells.Item(6, 14).Value = 1.010946152266434

$ws.Cells.Item(7, 2).Value = 1.02
$ws.Cells.Item(7, 3).Value = 1.015628647528756
$ws.Cells.Item(7, 4).Value = 1.021802555289802
$ws.Cells.Item(7, 5).Value = 1.017232158014215
$ws.Cells.Item(7, 9).Value = 1.026043130307829
$ws.Cells.Item(7, 10).Value = 1.020257992182585
$ws.Cells.Item(7, 11).Value = 1.02432226044508
$ws.Cells.Item(7, 12).Value = 1.019763824732977
$ws.Cells.Item(7, 14).Value = 1.010900587854566

$ws.Cells.Item(8, 2).Value = 1.02
$ws.Cells.Item(8, 3).Value = 1.014567360830581
$ws.Cells.Item(8, 4).Value = 1.021028701805798
$ws.Cells.Item(8, 5).Value = 1.016337816463709
$ws.Cells.Item(8, 9).Value = 1.025911497458653
$ws.Cells.Item(8, 10).Value = 1.01967628761025
$ws.Cells.Item(8, 11).Value = 1.023802721442219
$ws.Cells.Item(8, 12).Value = 1.01912543307562
$ws.Cells.Item(8, 14).Value = 1.010709766782213

$ws.Cells.Item(9, 2).Value = 1.02
$ws.Cells.Item(9, 3).Value = 1.012698020648709
$ws.Cells.Item(9, 4).Value = 1.019660828785941
$ws.Cells.Item(9, 5).Value = 1.014765807281991
$ws.Cells.Item(9, 9).Value = 1.025663670317693
$ws.Cells.Item(9, 10).Value = 1.018647159979196
$ws.Cells.Item(9, 11).Value = 1.022877644798702
$ws.Cells.Item(9, 12).Value = 1.017999246084912
$ws.Cells.Item(9, 14).Value = 1.010372154430315

$ws.Cells.Item(10, 2).Value = 1.02
$ws.Cells.Item(10, 3).Value = 1.011452420350442
$ws.Cells.Item(10, 4).Value = 1.018746181951012
$ws.Cells.Item(10, 5).Value = 1.01372055861418
$ws.Cells.Item(10, 9).Value = 1.025487842272465
$ws.Cells.Item(10, 10).Value = 1.017958407364963
$ws.Cells.Item(10, 11).Value = 1.022254580652733
$ws.Cells.Item(10, 12).Value = 1.017247694387752
$ws.Cells.Item(10, 14).Value = 1.010146191787924

$ws.Cells.Item(11, 2).Value = 1.02
$ws.Cells.Item(11, 3).Value = 1.010913231122607
$ws.Cells.Item(11, 4).Value = 1.018349507452007
$ws.Cells.Item(11, 5).Value = 1.013268631914525
$ws.Cells.Item(11, 9).Value = 1.025409203184091
$ws.Cells.Item(11, 10).Value = 1.017659554340371
$ws.Cells.Item(11, 11).Value = 1.021983301811832
$ws.Cells.Item(11, 12).Value = 1.016922101919302
$ws.Cells.Item(11, 14).Value = 1.010048142645139

$ws.Cells.Item(12, 2).Value = 1.02
$ws.Cells.Item(12, 3).Value = 1.010712978512408
$ws.Cells.Item(12, 4).Value = 1.018202072813729
$ws.Cells.Item(12, 5).Value = 1.01310086919875
$ws.Cells.Item(12, 9).Value = 1.025379617873677
$ws.Cells.Item(12, 10).Value = 1.017548455445025
$ws.Cells.Item(12, 11).Value = 1.021882314811298
$ws.Cells.Item(12, 12).Value = 1.016801138939023
$ws.Cells.Item(12, 14).Value = 1.01001169236006

$ws.Cells.Item(13, 2).Value = 1.02
$ws.Cells.Item(13, 3).Value = 1.010755932161591
$ws.Cells.Item(13, 4).Value = 1.018233702177762
$ws.Cells.Item(13, 5).Value = 1.013136850180327
$ws.Cells.Item(13, 9).Value = 1.025385980988091
$ws.Cells.Item(13, 10).Value = 1.017572290657422
$ws.Cells.Item(13, 11).Value = 1.021903986884523
$ws.Cells.Item(13, 12).Value = 1.016827086946565
$ws.Cells.Item(13, 14).Value = 1.010019512441392

$ws.Cells.Item(14, 2).Value = 1.02
$ws.Cells.Item(14, 3).Value = 1.010896677633831
$ws.Cells.Item(14, 4).Value = 1.018337322325217
$ws.Cells.Item(14, 5).Value = 1.01325476248439
$ws.Cells.Item(14, 9).Value = 1.025406765302347
$ws.Cells.Item(14, 10).Value = 1.017650372731293
$ws.Cells.Item(14, 11).Value = 1.021974958713648
$ws.Cells.Item(14, 12).Value = 1.016912103551482
$ws.Cells.Item(14, 14).Value = 1.010045130272114

$ws.Cells.Item(15, 2).Value = 1.02
$ws.Cells.Item(15, 3).Value = 1.010983399095573
$ws.Cells.Item(15, 4).Value = 1.018401153978898
$ws.Cells.Item(15, 5).Value = 1.013327425849108
$ws.Cells.Item(15, 9).Value = 1.025419521509727
$ws.Cells.Item(15, 10).Value = 1.017698469579974
$ws.Cells.Item(15, 11).Value = 1.022018657435329
$ws.Cells.Item(15, 12).Value = 1.016964482024593
$ws.Cells.Item(15, 14).Value = 1.010060910239411

$ws.Cells.Item(16, 2).Value = 1.02
$ws.Cells.Item(16, 3).Value = 1.011488208158249
$ws.Cells.Item(16, 4).Value = 1.018772494908339
$ws.Cells.Item(16, 5).Value = 1.013750565839126
$ws.Cells.Item(16, 9).Value = 1.02549300861116
$ws.Cells.Item(16, 10).Value = 1.01797822833924
$ws.Cells.Item(16, 11).Value = 1.022272553322307
$ws.Cells.Item(16, 12).Value = 1.017269299521391
$ws.Cells.Item(16, 14).Value = 1.010152694689005

$ws.Cells.Item(17, 2).Value = 1.02
$ws.Cells.Item(17, 3).Value = 1.011804906598835
$ws.Cells.Item(17, 4).Value = 1.019005261006382
$ws.Cells.Item(17, 5).Value = 1.014016171860491
$ws.Cells.Item(17, 9).Value = 1.025538435345091
$ws.Cells.Item(17, 10).Value = 1.0181535489656
$ws.Cells.Item(17, 11).Value = 1.022431418261652
$ws.Cells.Item(17, 12).Value = 1.017460460162707
$ws.Cells.Item(17, 14).Value = 1.010210213865862

$ws.Cells.Item(18, 2).Value = 1.02
$ws.Cells.Item(18, 3).Value = 1.011989647146872
$ws.Cells.Item(18, 4).Value = 1.019140968960345
$ws.Cells.Item(18, 5).Value = 1.014171160155064
$ws.Cells.Item(18, 9).Value = 1.025564690286827
$ws.Cells.Item(18, 10).Value = 1.01825575074422
$ws.Cells.Item(18, 11).Value = 1.02252393782118
$ws.Cells.Item(18, 12).Value = 1.017571944853989
$ws.Cells.Item(18, 14).Value = 1.01024374393882

$ws.Cells.Item(19, 2).Value = 1.02
$ws.Cells.Item(19, 3).Value = 1.012052641513334
$ws.Cells.Item(19, 4).Value = 1.019187231529152
$ws.Cells.Item(19, 5).Value = 1.014224018091104
$ws.Cells.Item(19, 9).Value = 1.02557360151567
$ws.Cells.Item(19, 10).Value = 1.018290588728908
$ws.Cells.Item(19, 11).Value = 1.022555460182337
$ws.Cells.Item(19, 12).Value = 1.01760995548181
$ws.Cells.Item(19, 14).Value = 1.010255173439659

$ws.Cells.Item(20, 2).Value = 1.02
$ws.Cells.Item(20, 3).Value = 1.011770926216896
$ws.Cells.Item(20, 4).Value = 1.01898029365751
$ws.Cells.Item(20, 5).Value = 1.013987668132235
$ws.Cells.Item(20, 9).Value = 1.025533586479708
$ws.Cells.Item(20, 10).Value = 1.018134744894714
$ws.Cells.Item(20, 11).Value = 1.022414388406283
$ws.Cells.Item(20, 12).Value = 1.017439952087456
$ws.Cells.Item(20, 14).Value = 1.010204044656533

$ws.Cells.Item(21, 2).Value = 1.02
$ws.Cells.Item(21, 3).Value = 1.010855230861673
$ws.Cells.Item(21, 4).Value = 1.018306811299658
$ws.Cells.Item(21, 5).Value = 1.013220037391769
$ws.Cells.Item(21, 9).Value = 1.025400655190597
$ws.Cells.Item(21, 10).Value = 1.017627382023722
$ws.Cells.Item(21, 11).Value = 1.021954065393395
$ws.Cells.Item(21, 12).Value = 1.016887068915844
$ws.Cells.Item(21, 14).Value = 1.010037587296891

$ws.Cells.Item(22, 2).Value = 1.02
$ws.Cells.Item(22, 3).Value = 1.01027964968742
$ws.Cells.Item(22, 4).Value = 1.017882834536591
$ws.Cells.Item(22, 5).Value = 1.012737994232426
$ws.Cells.Item(22, 9).Value = 1.025314905413665
$ws.Cells.Item(22, 10).Value = 1.017307853689626
$ws.Cells.Item(22, 11).Value = 1.021663358628423
$ws.Cells.Item(22, 12).Value = 1.016539314410843
$ws.Cells.Item(22, 14).Value = 1.009932752894315

$ws.Cells.Item(23, 2).Value = 1.02
$ws.Cells.Item(23, 3).Value = 1.010584761121711
$ws.Cells.Item(23, 4).Value = 1.018107642359834
$ws.Cells.Item(23, 5).Value = 1.012993477237327
$ws.Cells.Item(23, 9).Value = 1.025360568407746
$ws.Cells.Item(23, 10).Value = 1.017477291410924
$ws.Cells.Item(23, 11).Value = 1.021817588910869
$ws.Cells.Item(23, 12).Value = 1.016723677962095
$ws.Cells.Item(23, 14).Value = 1.009988344131177

$ws.Cells.Item(24, 2).Value = 1.02
$ws.Cells.Item(24, 3).Value = 1.011786280440464
$ws.Cells.Item(24, 4).Value = 1.018991575515147
$ws.Cells.Item(24, 5).Value = 1.014000547540121
$ws.Cells.Item(24, 9).Value = 1.025535778220178
$ws.Cells.Item(24, 10).Value = 1.018143241829781
$ws.Cells.Item(24, 11).Value = 1.022422083908971
$ws.Cells.Item(24, 12).Value = 1.017449218853739
$ws.Cells.Item(24, 14).Value = 1.010206832318184

$ws.Cells.Item(25, 2).Value = 1.02
$ws.Cells.Item(25, 3).Value = 1.013181186791001
$ws.Cells.Item(25, 4).Value = 1.020014947080753
$ws.Cells.Item(25, 5).Value = 1.015171731469196
$ws.Cells.Item(25, 9).Value = 1.025729614177676
$ws.Cells.Item(25, 10).Value = 1.018913690299751
$ws.Cells.Item(25, 11).Value = 1.023117924887494
$ws.Cells.Item(25, 12).Value = 1.018290532726119
$ws.Cells.Item(25, 14).Value = 1.010459593787363
